$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the range C2:C81 (mirrors the author selecting this range in the UI
# before clearing it) so the saved sheetView reflects the same selection.
$range = $ws.Range("C2:C81")
$range.Select()

# Clear the cell contents (equivalent to pressing the Delete key on the
# selection) - this empties the "PID" column values for rows 2-81 while
# leaving each cell's existing style/formatting untouched.
$range.ClearContents()
